{"js": "// CIV-17609: updated GA documents to display main claim number.\n// The \"claim details\" panel on the GA letters was labelled \"Claim number:\"\n// even though the merge field underneath it (<<caseNumber>>) resolves the\n// *case* number, not the claim number. Relabel it to \"Case number:\" so the\n// displayed text matches the value that actually gets merged in.\n//\n// Only the specific \"Claim number: <<caseNumber>>\" line should change -\n// the document separately contains an (unrelated, already-correct)\n// \"Case number: <<claimNumber>>\" line elsewhere that must be left alone,\n// so we search for the exact, case-sensitive \"Claim number: \" label text.\nconst body = context.document.body;\nconst matches = body.search(\"Claim number: \", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  // Replacing only the matched text (not the whole paragraph/run) keeps the\n  // existing character formatting and leaves the merge field untouched.\n  matches.items[i].insertText(\"Case number: \", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# CIV-17609: updated GA documents to display main claim number.\n# The \"claim details\" panel on the GA letters was labelled \"Claim number:\"\n# even though the merge field underneath it (<<caseNumber>>) resolves the\n# *case* number, not the claim number. Relabel it to \"Case number:\" so the\n# displayed text matches the value that actually gets merged in.\n#\n# Only the specific \"Claim number: <<caseNumber>>\" line should change - the\n# document separately contains an (unrelated, already-correct)\n# \"Case number: <<claimNumber>>\" line elsewhere that must be left alone, so\n# we search for the exact, case-sensitive \"Claim number: \" label text and\n# replace just that text (leaving the merge field and formatting intact).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Claim number: \"\n$find.Replacement.Text = \"Case number: \"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
